$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(''Chandra, Torch of Defiance'', [''{2}{R}{R}'', ''Legendary Planeswalker — Chandra'', ''+1: Exile the top card of your library. You may cast that card. If you don’t, Chandra, Torch of Defiance deals 2 damage to each opponent.'', ''+1: Add {R}{R}.'', ''−3: Chandra, Torch of Defiance deals 4 damage to target creature.'', ''−7: You get an emblem with “Whenever you cast a spell, this emblem deals 5 damage to any target.”'', ''Loyalty: 4''])'
$ws.Range("A3").Value = '(''Gideon of the Trials'', [''{1}{W}{W}'', ''Legendary Planeswalker — Gideon'', ''+1: Until your next turn, prevent all damage target permanent would deal.'', ''0: Until end of turn, Gideon of the Trials becomes a 4/4 Human Soldier creature with indestructible that’s still a planeswalker. Prevent all damage that would be dealt to him this turn.'', ''0: You get an emblem with “As long as you control a Gideon planeswalker, you can’t lose the game and your opponents can’t win the game.”'', ''Loyalty: 3''])'
$ws.Range("A4").Value = '(''Jace, Cunning Castaway'', [''{1}{U}{U}'', ''Legendary Planeswalker — Jace'', ''+1: Whenever one or more creatures you control deal combat damage to a player this turn, draw a card, then discard a card.'', ''−2: Create a 2/2 blue Illusion creature token with “When this creature becomes the target of a spell, sacrifice it.”'', ''−5: Create two tokens that are copies of Jace, Cunning Castaway, except they’re not legendary.'', ''Loyalty: 3''])'
$ws.Range("A5").Value = '(''Liliana, Untouched by Death'', [''{2}{B}{B}'', ''Legendary Planeswalker — Liliana'', ''+1: Mill three cards. If at least one Zombie card is milled this way, each opponent loses 2 life and you gain 2 life.'', ''−2: Target creature gets -X/-X until end of turn, where X is the number of Zombies you control.'', ''−3: You may cast Zombie spells from your graveyard this turn.'', ''Loyalty: 4''])'
$ws.Range("A6").Value = '(''Nissa, Vital Force'', [''{3}{G}{G}'', ''Legendary Planeswalker — Nissa'', ''+1: Untap target land you control. Until your next turn, it becomes a 5/5 Elemental creature with haste. It’s still a land.'', ''−3: Return target permanent card from your graveyard to your hand.'', ''−6: You get an emblem with “Whenever a land enters the battlefield under your control, you may draw a card.”'', ''Loyalty: 5''])'

$ws.Range("A7:A37").Clear()
